# Insert a new row above row 277, shifting existing rows 277:335 down to 278:336.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(277).Insert()

# Populate the newly inserted row 277 with the new record.
$ws.Range("A277").Value = 5
$ws.Range("B277").Value = "Macroferia Regional de Talca"
$ws.Range("C277").Value = "Maule"
$ws.Range("D277").Value = 44711
$ws.Range("E277").Value = 7
$ws.Range("F277").Value = 100112023
$ws.Range("G277").Value = "Brócoli"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 3000
$ws.Range("K277").Value = 1000
$ws.Range("L277").Value = 1000
$ws.Range("M277").Value = 1000
$ws.Range("N277").Value = "$/unidad"
$ws.Range("O277").Value = "Región del Maule"
$ws.Range("P277").Value = 1000
$ws.Range("Q277").Value = 1
$ws.Range("R277").Value = "Hortaliza"
